$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: mark the "computational thinking,scratch" CS50 entry as complete,
# recording the completion date/time in the two new columns F/G.
$ws.Range("F1").Value = "complete on date"
$ws.Range("G1").Value = "complete on time"

$ws.Range("F2").NumberFormat = "m/d/yy"
$ws.Range("F2").Formula = "=TODAY()"
$ws.Range("G2").NumberFormat = "h:mm AM/PM"
$ws.Range("G2").Value = 0.41666666666666669

# Row 3: a new entry for the C program that was started.
$ws.Range("A3").Value = "CS50"
$ws.Range("B3").Value = "writing program in c"
$ws.Range("C3").NumberFormat = "m/d/yy"
$ws.Range("C3").Formula = "=TODAY()"
$ws.Range("D3").NumberFormat = "h:mm AM/PM"
$ws.Range("D3").Value = 0.46597222222222223
$ws.Range("E3").Value = "youtube"

$ws.Columns("F:F").ColumnWidth = $ws.Columns("C:C").ColumnWidth

$ws.Range("D4").Select() | Out-Null
